$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.345130039178628
$ws.Range("C2").Value = 0.2415607546647607
$ws.Range("D2").Value = 0.5752131711939796
$ws.Range("E2").Value = 0.2038491611598161
$ws.Range("G2").Value = 2.630987599011092
$ws.Range("H2").Value = 2.02042662921491
$ws.Range("J2").Value = 0.08525212067916854
$ws.Range("M2").Value = 0.5991978372868729

$ws.Range("B3").Value = 1.268847915797494
$ws.Range("C3").Value = 0.2234163728232943
$ws.Range("D3").Value = 0.5714660739660786
$ws.Range("E3").Value = 0.2032244090067898
$ws.Range("G3").Value = 2.566958220720636
$ws.Range("H3").Value = 1.99676484978562
$ws.Range("J3").Value = 0.08535194105846244
$ws.Range("M3").Value = 0.5794400758944036

$ws.Range("B4").Value = 1.222965749027509
$ws.Range("C4").Value = 0.2124525988794801
$ws.Range("D4").Value = 0.5694355102835118
$ws.Range("E4").Value = 0.2029249450473714
$ws.Range("G4").Value = 2.528969958683973
$ws.Range("H4").Value = 1.983123237838754
$ws.Range("J4").Value = 0.08543874495920178
$ws.Range("M4").Value = 0.5677167756999211

$ws.Range("B5").Value = 1.204507760901976
$ws.Range("C5").Value = 0.2080290173423123
$ws.Range("D5").Value = 0.5686758551162825
$ws.Range("E5").Value = 0.2028240339722593
$ws.Range("G5").Value = 2.513819723511489
$ws.Range("H5").Value = 1.97778583440072
$ws.Range("J5").Value = 0.08548052533149786
$ws.Range("M5").Value = 0.5630417536437804

$ws.Range("B6").Value = 1.201457258739396
$ws.Range("C6").Value = 0.2072971502211658
$ws.Range("D6").Value = 0.5685538074461505
$ws.Range("E6").Value = 0.2028085525696532
$ws.Range("G6").Value = 2.511323895416382
$ws.Range("H6").Value = 1.976912911810871
$ws.Range("J6").Value = 0.08548784965605094
$ws.Range("M6").Value = 0.5622716413752684

$ws.Range("B7").Value = 1.222715850084711
$ws.Range("C7").Value = 0.2123927620864379
$ws.Range("D7").Value = 0.5694249908778062
$ws.Range("E7").Value = 0.2029234986432442
$ws.Range("G7").Value = 2.528764304154663
$ws.Range("H7").Value = 1.983050360031285
$ws.Range("J7").Value = 0.08543928249299526
$ws.Range("M7").Value = 0.5676533128109824

$ws.Range("B8").Value = 1.318629116238867
$ws.Range("C8").Value = 0.2352676860218423
$ws.Range("D8").Value = 0.5738650225309669
$ws.Range("E8").Value = 0.2036162628363094
$ws.Range("G8").Value = 2.608633125672753
$ws.Range("H8").Value = 2.012083173516942
$ws.Range("J8").Value = 0.08528123804657994
$ws.Range("M8").Value = 0.5923004864658736

$ws.Range("B9").Value = 1.514347671687062
$ws.Range("C9").Value = 0.2815437472201836
$ws.Range("D9").Value = 0.5847222714084523
$ws.Range("E9").Value = 0.2056442016785809
$ws.Range("G9").Value = 2.775929471921415
$ws.Range("H9").Value = 2.076117843628879
$ws.Range("J9").Value = 0.08517419653825087
$ws.Range("M9").Value = 0.6438881113305825

$ws.Range("B10").Value = 1.662887505595165
$ws.Range("C10").Value = 0.3164318825472492
$ws.Range("D10").Value = 0.5940211113299085
$ws.Range("E10").Value = 0.2075452502255004
$ws.Range("G10").Value = 2.905574761223534
$ws.Range("H10").Value = 2.127590167287593
$ws.Range("J10").Value = 0.08521993130423056
$ws.Range("M10").Value = 0.6838027715396322

$ws.Range("B11").Value = 1.731514288898609
$ws.Range("C11").Value = 0.3325021846420952
$ws.Range("D11").Value = 0.5985410422225357
$ws.Range("E11").Value = 0.2085000580349821
$ws.Range("G11").Value = 2.966065589227412
$ws.Range("H11").Value = 2.151988844281163
$ws.Range("J11").Value = 0.08526790447059085
$ws.Range("M11").Value = 0.7024049873371894

$ws.Range("B12").Value = 1.757654679647942
$ws.Range("C12").Value = 0.3386166798336774
$ws.Range("D12").Value = 0.6002944872197418
$ws.Range("E12").Value = 0.2088746137108508
$ws.Range("G12").Value = 2.989193478030813
$ws.Range("H12").Value = 2.161371073852592
$ws.Range("J12").Value = 0.08528998970091806
$ws.Range("M12").Value = 0.7095136024352229

$ws.Range("B13").Value = 1.752018058572901
$ws.Range("C13").Value = 0.3372985177869907
$ws.Range("D13").Value = 0.5999149876827232
$ws.Range("E13").Value = 0.2087933679144172
$ws.Range("G13").Value = 2.984202573157347
$ws.Range("H13").Value = 2.159344059174714
$ws.Range("J13").Value = 0.08528505876458325
$ws.Range("M13").Value = 0.7079797661514959

$ws.Range("B14").Value = 1.733661803003827
$ws.Range("C14").Value = 0.333004644283335
$ws.Range("D14").Value = 0.5986844595844332
$ws.Range("E14").Value = 0.2085306123280439
$ws.Range("G14").Value = 2.967963878483999
$ws.Range("H14").Value = 2.152757851584909
$ws.Range("J14").Value = 0.08526964282619431
$ws.Range("M14").Value = 0.7029885252110262

$ws.Range("B15").Value = 1.722438016387684
$ws.Range("C15").Value = 0.3303783121557444
$ws.Range("D15").Value = 0.597936180110878
$ws.Range("E15").Value = 0.208371360138738
$ws.Range("G15").Value = 2.958046137769713
$ws.Range("H15").Value = 2.148742275252459
$ws.Range("J15").Value = 0.08526071083033315
$ws.Range("M15").Value = 0.6999396394668338

$ws.Range("B16").Value = 1.658423917583605
$ws.Range("C16").Value = 0.3153856932099472
$ws.Range("D16").Value = 0.5937315696760663
$ws.Range("E16").Value = 0.2074846665464847
$ws.Range("G16").Value = 2.901652321165813
$ws.Range("H16").Value = 2.126015580560619
$ws.Range("J16").Value = 0.08521734386420121
$ws.Range("M16").Value = 0.682596057461339

$ws.Range("B17").Value = 1.619424438497845
$ws.Range("C17").Value = 0.3062395385485388
$ws.Range("D17").Value = 0.5912265295273755
$ws.Range("E17").Value = 0.2069637959123654
$ws.Range("G17").Value = 2.867446954991863
$ws.Range("H17").Value = 2.112326511702094
$ws.Range("J17").Value = 0.08519770655919601
$ws.Range("M17").Value = 0.6720705395870539

$ws.Range("B18").Value = 1.597092223384379
$ws.Range("C18").Value = 0.3009976709352316
$ws.Range("D18").Value = 0.5898129677242707
$ws.Range("E18").Value = 0.2066726742841922
$ws.Range("G18").Value = 2.847915274695453
$ws.Range("H18").Value = 2.10454541534483
$ws.Range("J18").Value = 0.08518896815910892
$ws.Range("M18").Value = 0.6660583858080074

$ws.Range("B19").Value = 1.589547931108712
$ws.Range("C19").Value = 0.2992260765874164
$ws.Range("D19").Value = 0.589339038829138
$ws.Range("E19").Value = 0.2065755587134319
$ws.Range("G19").Value = 2.841326523485776
$ws.Range("H19").Value = 2.101926709375107
$ws.Range("J19").Value = 0.08518644819780974
$ws.Range("M19").Value = 0.6640299491643304

$ws.Range("B20").Value = 1.623565718912459
$ws.Range("C20").Value = 0.3072112185259073
$ws.Range("D20").Value = 0.5914903716912079
$ws.Range("E20").Value = 0.2070183665945358
$ws.Range("G20").Value = 2.87107341729785
$ws.Range("H20").Value = 2.113774152645647
$ws.Range("J20").Value = 0.08519953231354194
$ws.Range("M20").Value = 0.6731866652656748

$ws.Range("B21").Value = 1.739049320516528
$ws.Range("C21").Value = 0.3342650687052355
$ws.Range("D21").Value = 0.5990447584947844
$ws.Range("E21").Value = 0.208607437066771
$ws.Range("G21").Value = 2.972727541671304
$ws.Range("H21").Value = 2.154688487668921
$ws.Range("J21").Value = 0.08527406440721208
$ws.Range("M21").Value = 0.7044528237546785

$ws.Range("B22").Value = 1.815416427913021
$ws.Range("C22").Value = 0.3521156408045272
$ws.Range("D22").Value = 0.6042259742064573
$ws.Range("E22").Value = 0.2097217290180922
$ws.Range("G22").Value = 3.040456108405806
$ws.Range("H22").Value = 2.182262492350105
$ws.Range("J22").Value = 0.08534562322950023
$ws.Range("M22").Value = 0.7252623966082723

$ws.Range("B23").Value = 1.774575856935598
$ws.Range("C23").Value = 0.3425728492595397
$ws.Range("D23").Value = 0.6014382806687308
$ws.Range("E23").Value = 0.2091200640288164
$ws.Range("G23").Value = 3.004188726868335
$ws.Range("H23").Value = 2.167468894088074
$ws.Range("J23").Value = 0.08530533629178194
$ws.Range("M23").Value = 0.7141214593590774

$ws.Range("B24").Value = 1.621693168582624
$ws.Range("C24").Value = 0.3067718710161671
$ws.Range("D24").Value = 0.5913710057484991
$ws.Range("E24").Value = 0.2069936692347305
$ws.Range("G24").Value = 2.86943347831928
$ws.Range("H24").Value = 2.11311939728833
$ws.Range("J24").Value = 0.08519869894369236
$ws.Range("M24").Value = 0.6726819430673316

$ws.Range("B25").Value = 1.460574103222655
$ws.Range("C25").Value = 0.2688705193271801
$ws.Range("D25").Value = 0.5815537574807053
$ws.Range("E25").Value = 0.2050236160485213
$ws.Range("G25").Value = 2.7295054676012
$ws.Range("H25").Value = 2.058024563387363
$ws.Range("J25").Value = 0.08518136322630454
$ws.Range("M25").Value = 0.6295809060016779
